$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1802.5
$ws.Range("I41").Value = 1691.75
$ws.Range("J41").Value = 2024
$ws.Range("K41").Value = 1691.75
$ws.Range("L41").Value = 2024
$ws.Range("M41").Value = -1251.75
$ws.Range("N41").Value = -2904

$ws.Range("H69").Value = 9749.5
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9749.5
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 29248.5
$ws.Range("N69").Value = -30996.5
$ws.Range("M69").ClearContents()

$ws.Range("H70").Value = 4066533.8
$ws.Range("I70").Value = 12195121
$ws.Range("J70").Value = 2240
$ws.Range("K70").Value = 36585363
$ws.Range("L70").Value = 6720
$ws.Range("M70").Value = -36585093
$ws.Range("N70").Value = -7260

$ws.Range("H72").Value = 9749.5
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9749.5
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 87745.5
$ws.Range("N72").Value = -96481.5
$ws.Range("M72").ClearContents()

$ws.Range("H73").Value = 4066533.8
$ws.Range("I73").Value = 12195121
$ws.Range("J73").Value = 2240
$ws.Range("K73").Value = 36585363
$ws.Range("L73").Value = 6720
$ws.Range("M73").Value = -36584427
$ws.Range("N73").Value = -8592

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H118").Value = 1354
$ws.Range("J118").Value = 4000
$ws.Range("L118").Value = 12000
$ws.Range("N118").Value = -15314

$ws.Range("H135").Value = 4959.5
$ws.Range("I135").Value = 360.25
$ws.Range("K135").Value = 3242.25
$ws.Range("M135").Value = -707.25

$ws.Range("H138").Value = 2626.373
$ws.Range("I138").Value = 1643.7428
$ws.Range("K138").Value = 4931.2284
$ws.Range("M138").Value = 208.7716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H103").Value = 69964
$ws.Range("J103").Value = 69964
$ws.Range("L103").Value = 69964
$ws.Range("N103").Value = -72308

$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 2565715.2
$ws.Range("I132").Value = 1458.4828
$ws.Range("J132").Value = 10002060
$ws.Range("K132").Value = 4375.4484
$ws.Range("L132").Value = 30006180
$ws.Range("M132").Value = -1845.4484
$ws.Range("N132").Value = -30011240

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2635
$ws.Range("I86").Value = 1825.8
$ws.Range("K86").Value = 1825.8
$ws.Range("M86").Value = -702.8

$ws.Range("H89").Value = 2635
$ws.Range("I89").Value = 1825.8
$ws.Range("K89").Value = 9129
$ws.Range("M89").Value = -3513

$ws.Range("H134").Value = 2175969.8
$ws.Range("I134").Value = 1972.225
$ws.Range("J134").Value = 16669286
$ws.Range("K134").Value = 5916.674999999999
$ws.Range("L134").Value = 50007858
$ws.Range("M134").Value = -3381.674999999999
$ws.Range("N134").Value = -50012928

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20410818
$ws.Range("I31").Value = 28573612
$ws.Range("J31").Value = 3831.2856
$ws.Range("K31").Value = 28573612
$ws.Range("L31").Value = 3831.2856
$ws.Range("M31").Value = -28573317
$ws.Range("N31").Value = -4421.2856

$ws.Range("H34").Value = 20410818
$ws.Range("I34").Value = 28573612
$ws.Range("J34").Value = 3831.2856
$ws.Range("K34").Value = 28573612
$ws.Range("L34").Value = 3831.2856
$ws.Range("M34").Value = -28573410
$ws.Range("N34").Value = -4235.2856

$ws.Range("H86").Value = 7915.5
$ws.Range("I86").Value = 6349.5
$ws.Range("K86").Value = 6349.5
$ws.Range("M86").Value = -5226.5

$ws.Range("H89").Value = 7915.5
$ws.Range("I89").Value = 6349.5
$ws.Range("K89").Value = 31747.5
$ws.Range("M89").Value = -26131.5

$ws.Range("H107").Value = 1317.0834
$ws.Range("I107").Value = 440.53333
$ws.Range("J107").Value = 2778
$ws.Range("K107").Value = 440.53333
$ws.Range("L107").Value = 2778
$ws.Range("M107").Value = 1479.46667
$ws.Range("N107").Value = -6618

$ws.Range("H132").Value = 2394.2
$ws.Range("I132").Value = 2110.8823
$ws.Range("K132").Value = 6332.646900000001
$ws.Range("M132").Value = -3802.646900000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4929.0835
$ws.Range("I131").Value = 2239.7273
$ws.Range("K131").Value = 6719.1819
$ws.Range("M131").Value = -1679.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2229555.2
$ws.Range("I24").Value = 2504499.8
$ws.Range("K24").Value = 2504499.8
$ws.Range("M24").Value = -2504326.8

$ws.Range("H102").Value = 1903.75
$ws.Range("I102").Value = 1920.6666
$ws.Range("J102").Value = 1650
$ws.Range("K102").Value = 1920.6666
$ws.Range("L102").Value = 1650
$ws.Range("M102").Value = -298.6666
$ws.Range("N102").Value = -4894

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 1672666.6
$ws.Range("I20").Value = 1672666.6
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1672666.6
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1672440.6
$ws.Range("N20").ClearContents()

$ws.Range("H82").Value = 4687.4546
$ws.Range("I82").Value = 1307.8334
$ws.Range("J82").Value = 8743
$ws.Range("K82").Value = 1307.8334
$ws.Range("L82").Value = 8743
$ws.Range("M82").Value = -946.8334
$ws.Range("N82").Value = -9465

$ws.Range("H85").Value = 4687.4546
$ws.Range("I85").Value = 1307.8334
$ws.Range("J85").Value = 8743
$ws.Range("K85").Value = 1307.8334
$ws.Range("L85").Value = 8743
$ws.Range("M85").Value = -59.83339999999998
$ws.Range("N85").Value = -11239

$ws.Range("H132").Value = 3211.717
$ws.Range("I132").Value = 1681.4688
$ws.Range("J132").Value = 5543.524
$ws.Range("K132").Value = 5044.4064
$ws.Range("L132").Value = 16630.572
$ws.Range("M132").Value = -2514.4064
$ws.Range("N132").Value = -21690.572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8500
$ws.Range("J14").Value = 10000
$ws.Range("L14").Value = 10000
$ws.Range("N14").Value = -10336

$ws.Range("H15").Value = 9000
$ws.Range("J15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("N15").Value = -9576

$ws.Range("H37").Value = 19999.666
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 39999
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 39999
$ws.Range("M37").Value = -9797
$ws.Range("N37").Value = -40405

$ws.Range("H52").Value = 30042
$ws.Range("I52").Value = 30042
$ws.Range("K52").Value = 30042
$ws.Range("M52").Value = -29816

$ws.Range("H58").Value = 64606.5
$ws.Range("J58").Value = 64606.5
$ws.Range("L58").Value = 64606.5
$ws.Range("N58").Value = -65222.5

$ws.Range("H62").Value = 12935.429
$ws.Range("I62").Value = 4799
$ws.Range("K62").Value = 4799
$ws.Range("M62").Value = -4175

$ws.Range("H65").Value = 12935.429
$ws.Range("I65").Value = 4799
$ws.Range("K65").Value = 23995
$ws.Range("M65").Value = -20875

$ws.Range("H81").Value = 5872.5
$ws.Range("I81").Value = 1745
$ws.Range("K81").Value = 3490
$ws.Range("M81").Value = -2429

$ws.Range("H84").Value = 5872.5
$ws.Range("I84").Value = 1745
$ws.Range("K84").Value = 17450
$ws.Range("M84").Value = -12146
